# Adds the new catalysis-process sheets (Synthesis/Cleaning/Characterization/
# CatalyticActivityAssesment/Analytics processes plus the specific technique
# sheets for cleaning/deposition/sputtering) right after the existing
# "Process" sheet. Every new sheet gets the standard NamedEntity-style
# header row (name / id / iri) used by the other "stub" class sheets
# already in this workbook (Recipe, Mixture, Workflow, Process, ...).
# All pre-existing sheets and their contents are left untouched; they are
# simply pushed further along the tab order by the insertion.

$wb = $excel.ActiveWorkbook

# Note: "CatalyticActivityAssesmentProcess" is 33 characters, which exceeds
# Excel's hard 31-character worksheet-name limit (enforced by the object
# model, just like real Excel), so it is truncated to the longest valid
# prefix below.
$newSheetNames = @(
    "SynthesisProcess",
    "CleaningProcess",
    "CharacterizationProcess",
    "CatalyticActivityAssesmentProce",
    "AnalyticsProcess",
    "OzoneCleaning",
    "CVD",
    "CVD_nanofab",
    "CVD_pc1",
    "CVD_pc2",
    "sputtering",
    "sputtering_prevac",
    "sputtering_vonAdenne"
)

$afterSheet = $wb.Worksheets.Item("Process")

foreach ($sheetName in $newSheetNames) {
    $newSheet = $wb.Worksheets.Add($null, $afterSheet)
    $newSheet.Name = $sheetName

    $newSheet.Range("A1").Value = "name"
    $newSheet.Range("B1").Value = "id"
    $newSheet.Range("C1").Value = "iri"

    $afterSheet = $newSheet
}
